# Spending Drivers - Higher - Instructional and Total Employment and Payroll
# Commit: "new data depot data files (removing box links from excel)"
#
# Changes applied:
#  1. Make "Employment" the active/selected sheet (was "Payroll").
#  2. Update the remembered cell selection on each sheet
#     (Employment -> E6, Payroll -> E7).
#  3. Give column A (the State/Acronym column, identical on both sheets) an
#     explicit best-fit width, as happens when a user auto-fits that column
#     after trimming it down (e.g. once the old hyperlinked "box.com" source
#     links were removed).

$wb = $excel.ActiveWorkbook

$wsEmployment = $wb.Worksheets.Item("Employment")
$wsPayroll    = $wb.Worksheets.Item("Payroll")

# --- Column A width (State column) on both sheets ---------------------
# Target stored width is ~18.71 characters (Excel's AutoFit result for
# "District of Columbia" in the default font). Set explicitly so both
# sheets match.
$wsEmployment.Columns("A:A").ColumnWidth = 17.8333333
$wsPayroll.Columns("A:A").ColumnWidth = 17.8333333

# --- Selection per sheet ------------------------------------------------
$wsEmployment.Range("E6").Select() | Out-Null
$wsPayroll.Range("E7").Select() | Out-Null

# --- Active sheet / tab --------------------------------------------------
# Employment becomes the active (visible-on-open) sheet instead of Payroll.
$wsEmployment.Activate() | Out-Null
$wsEmployment.Range("E6").Select() | Out-Null
